$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a handful of numeric-looking Price cells to remain Text,
# matching the original inlineStr storage (e.g. "1.00" must not become 1).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '29.904.10'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.634.61'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.86%  '
$ws.Range("D5").Value = '215.42'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("D8").Value = '28.71'
$ws.Range("E8").Value = '  -1.90%  '
$ws.Range("D9").Value = '0.261'
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("D10").Value = '0.0610'
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("D12").Value = '1.869.55'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '1.633.73'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '0.588'
$ws.Range("E14").Value = '  +3.97%  '
$ws.Range("D15").Value = '9.52'
$ws.Range("E15").Value = '  +6.43%  '
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").Value = '29.913.11'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = '65.00'
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").Value = '240.38'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '0.0₃0704'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").Value = '9.89'
$ws.Range("E22").Value = '  +2.66%  '
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("E24").Value = '  +2.90%  '
$ws.Range("D25").Value = '158.21'
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("D26").Value = '15.54'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("E27").Value = '  -1.14%  '
$ws.Range("D28").Value = '6.63'
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("E32").Value = '  +1.41%  '
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("D34").Value = '1.424.32'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("E35").Value = '  +4.22%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("E37").Value = '  -3.16%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").Value = '76.08'
$ws.Range("E40").Value = '  +9.54%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("E42").Value = '  -0.35%  '
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("D47").Value = '1.777.21'
$ws.Range("E47").Value = '  +0.96%  '
$ws.Range("D48").Value = '5.34'
$ws.Range("E48").Value = '  -1.93%  '
$ws.Range("D49").Value = '48.71'
$ws.Range("E49").Value = '  -9.27%  '
$ws.Range("D50").Value = '92.73'
$ws.Range("E50").Value = '  +5.15%  '
$ws.Range("E51").Value = '  +8.52%  '
